$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Commit: "unify the conception of DataNode, DataTable, Entity."
# The only intentional content-level change in the XML diff is the
# worksheet's display name ("Property1" -> "DataNode"); everything else
# in the diff (fileVersion/rupBuild, absPath, window geometry, xr/xr2/xr9
# revision GUIDs, tiny column-width deltas, the extra phonetic-guide font,
# cellStyle locale label, selection position, etc.) is Excel-version/
# environment re-save noise, not a deliberate edit.
$ws.Name = "DataNode"

